$wb = $excel.ActiveWorkbook

# Sheet ALC, row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 192
$ws.Range("J12").Value = 300
$ws.Range("L12").Value = 300
$ws.Range("N12").Value = -640

# Sheet ALC, row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 536.16327
$ws.Range("I15").Value = 536.16327
$ws.Range("K15").Value = 1608.48981
$ws.Range("M15").Value = -1439.48981

# Sheet ALC, row 110
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 67957.60000000001
$ws.Range("J110").Value = 67957.60000000001
$ws.Range("L110").Value = 67957.60000000001
$ws.Range("N110").Value = -76137.60000000001

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8321735
$ws.Range("I116").Value = 9968044
$ws.Range("K116").Value = 9968044
$ws.Range("M116").Value = -9964602

# Sheet ALC, row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 79269.125
$ws.Range("J123").Value = 79269.125
$ws.Range("L123").Value = 79269.125
$ws.Range("N123").Value = -89069.125

# Sheet ALC, row 136
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 77150.336
$ws.Range("J136").Value = 77150.336
$ws.Range("L136").Value = 77150.336
$ws.Range("N136").Value = -87350.336

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2728.875
$ws.Range("J138").Value = 3321.15
$ws.Range("L138").Value = 9963.450000000001
$ws.Range("N138").Value = -20243.45

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 25004542
$ws.Range("I45").Value = 5903.6665
$ws.Range("J45").Value = 62502500
$ws.Range("K45").Value = 5903.6665
$ws.Range("L45").Value = 62502500
$ws.Range("M45").Value = -5526.6665
$ws.Range("N45").Value = -62503254

# Sheet ARM, row 52
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 52897
$ws.Range("J52").Value = 52897
$ws.Range("L52").Value = 52897
$ws.Range("N52").Value = -53533

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 47280.953
$ws.Range("I61").Value = 894.1429000000001
$ws.Range("J61").Value = 128457.875
$ws.Range("K61").Value = 894.1429000000001
$ws.Range("L61").Value = 128457.875
$ws.Range("M61").Value = -682.1429000000001
$ws.Range("N61").Value = -128881.875

# Sheet ARM, row 104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 39997
$ws.Range("J104").Value = 39997
$ws.Range("L104").Value = 39997
$ws.Range("N104").Value = -46985

# Sheet ARM, row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1541.1
$ws.Range("I110").Value = 938.875
$ws.Range("K110").Value = 938.875
$ws.Range("M110").Value = 1106.125

# Sheet ARM, row 121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 51865
$ws.Range("J121").Value = 51865
$ws.Range("L121").Value = 51865
$ws.Range("N121").Value = -55359

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 47280.953
$ws.Range("I136").Value = 894.1429000000001
$ws.Range("J136").Value = 128457.875
$ws.Range("K136").Value = 2682.4287
$ws.Range("L136").Value = 385373.625
$ws.Range("M136").Value = -132.4287000000004
$ws.Range("N136").Value = -390473.625

# Sheet BSM, row 6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 5500
$ws.Range("I6").Value = 4500
$ws.Range("K6").Value = 4500
$ws.Range("M6").Value = -4387

# Sheet BSM, row 50
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 72675.25
$ws.Range("J50").Value = 79997.336
$ws.Range("L50").Value = 79997.336
$ws.Range("N50").Value = -81145.336

# Sheet BSM, row 52
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 99988
$ws.Range("J52").Value = 99988
$ws.Range("L52").Value = 99988
$ws.Range("N52").Value = -100514

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1748.5238
$ws.Range("I94").Value = 1485.45
$ws.Range("K94").Value = 1485.45
$ws.Range("M94").Value = -1034.45

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1426.6666
$ws.Range("I107").Value = 1200.6666
$ws.Range("J107").Value = 2330.6667
$ws.Range("K107").Value = 1200.6666
$ws.Range("L107").Value = 2330.6667
$ws.Range("M107").Value = 719.3334
$ws.Range("N107").Value = -6170.6667

# Sheet BSM, row 109
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 91992.8
$ws.Range("J109").Value = 91992.8
$ws.Range("L109").Value = 91992.8
$ws.Range("N109").Value = -94766.8

# Sheet BSM, row 110
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 72897
$ws.Range("J110").Value = 72897
$ws.Range("L110").Value = 72897
$ws.Range("N110").Value = -81077

# Sheet BSM, row 117
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 87347.5
$ws.Range("J117").Value = 87347.5
$ws.Range("L117").Value = 87347.5
$ws.Range("N117").Value = -96525.5

# Sheet BSM, row 119
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 50000.5
$ws.Range("J119").Value = 50000.5
$ws.Range("L119").Value = 50000.5
$ws.Range("N119").Value = -59676.5

# Sheet BSM, row 121
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H121").Value = 99988
$ws.Range("J121").Value = 99988
$ws.Range("L121").Value = 99988
$ws.Range("N121").Value = -103482

# Sheet BSM, row 122
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 76346.875
$ws.Range("J122").Value = 76346.875
$ws.Range("L122").Value = 76346.875
$ws.Range("N122").Value = -86146.875

# Sheet BSM, row 129
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

# Sheet CRP, row 9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 25511.5
$ws.Range("J9").Value = 25511.5
$ws.Range("L9").Value = 25511.5
$ws.Range("N9").Value = -25847.5

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1506
$ws.Range("I107").Value = 100
$ws.Range("J107").Value = 1740.3334
$ws.Range("K107").Value = 100
$ws.Range("L107").Value = 1740.3334
$ws.Range("M107").Value = 1820
$ws.Range("N107").Value = -5580.3334

# Sheet CRP, row 108
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 57766.9
$ws.Range("J108").Value = 61963.223
$ws.Range("L108").Value = 61963.223
$ws.Range("N108").Value = -69643.223

# Sheet CRP, row 117
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H117").Value = 42998
$ws.Range("J117").Value = 42998
$ws.Range("L117").Value = 42998
$ws.Range("N117").Value = -52176

# Sheet CRP, row 119
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 64150.8
$ws.Range("J119").Value = 64150.8
$ws.Range("L119").Value = 64150.8
$ws.Range("N119").Value = -73826.8

# Sheet CRP, row 125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 52365.2
$ws.Range("J125").Value = 52365.2
$ws.Range("L125").Value = 52365.2
$ws.Range("N125").Value = -57285.2

# Sheet CUL, row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 4999
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Sheet GSM, row 26
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 9225.777
$ws.Range("J26").Value = 9225.777
$ws.Range("L26").Value = 9225.777
$ws.Range("N26").Value = -9785.777

# Sheet GSM, row 50
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 9225.777
$ws.Range("J50").Value = 9225.777
$ws.Range("L50").Value = 9225.777
$ws.Range("N50").Value = -10221.777

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125003200
$ws.Range("I80").Value = 200002820
$ws.Range("K80").Value = 200002820
$ws.Range("M80").Value = -200001822

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 125003200
$ws.Range("I83").Value = 200002820
$ws.Range("K83").Value = 1000014100
$ws.Range("M83").Value = -1000009108

# Sheet GSM, row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1005.56665
$ws.Range("I97").Value = 572.8077
$ws.Range("J97").Value = 3818.5
$ws.Range("K97").Value = 572.8077
$ws.Range("L97").Value = 3818.5
$ws.Range("M97").Value = -76.80769999999995
$ws.Range("N97").Value = -4810.5

# Sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1133.5555
$ws.Range("I107").Value = 1033.6666
$ws.Range("J107").Value = 1333.3334
$ws.Range("K107").Value = 1033.6666
$ws.Range("L107").Value = 1333.3334
$ws.Range("M107").Value = 886.3334
$ws.Range("N107").Value = -5173.3334

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4320.5713
$ws.Range("I126").Value = 2231.3333
$ws.Range("J126").Value = 5887.5
$ws.Range("K126").Value = 6693.999899999999
$ws.Range("L126").Value = 17662.5
$ws.Range("M126").Value = -4223.999899999999
$ws.Range("N126").Value = -22602.5

# Sheet GSM, row 134
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 63329.668
$ws.Range("J134").Value = 63329.668
$ws.Range("L134").Value = 189989.004
$ws.Range("N134").Value = -195059.004

# Sheet GSM, row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 59661
$ws.Range("J140").Value = 98983
$ws.Range("L140").Value = 98983
$ws.Range("N140").Value = -109343

# Sheet LTW, row 118
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H118").Value = 55130.4
$ws.Range("J118").Value = 55130.4
$ws.Range("L118").Value = 55130.4
$ws.Range("N118").Value = -58444.4

# Sheet LTW, row 123
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 77832.664
$ws.Range("J123").Value = 77832.664
$ws.Range("L123").Value = 77832.664
$ws.Range("N123").Value = -87632.664

# Sheet LTW, row 125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# Sheet WVR, row 121
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 41595.668
$ws.Range("J121").Value = 46747.5
$ws.Range("L121").Value = 46747.5
$ws.Range("N121").Value = -50241.5

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2420.2856
$ws.Range("I122").Value = 2034.1904
$ws.Range("K122").Value = 6102.5712
$ws.Range("M122").Value = -3652.5712
